$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (G=27772)
$ws.Range("H28").Value = 1859.3
$ws.Range("I28").Value = 2070.9285
$ws.Range("J28").Value = 1365.5
$ws.Range("K28").Value = 2070.9285
$ws.Range("L28").Value = 1365.5
$ws.Range("M28").Value = -1585.9285
$ws.Range("N28").Value = -2335.5

# Row 62 (G=27781)
$ws.Range("H62").Value = 2388.3333
$ws.Range("I62").Value = 2070.7144
$ws.Range("K62").Value = 2070.7144
$ws.Range("M62").Value = -1446.7144

# Row 65 (G=27781)
$ws.Range("H65").Value = 2388.3333
$ws.Range("I65").Value = 2070.7144
$ws.Range("K65").Value = 10353.572
$ws.Range("M65").Value = -7233.572

# Row 98 (G=36237)
$ws.Range("H98").Value = 746.43335
$ws.Range("I98").Value = 768.0345
$ws.Range("J98").Value = 120
$ws.Range("K98").Value = 768.0345
$ws.Range("L98").Value = 120
$ws.Range("M98").Value = 729.9655
$ws.Range("N98").Value = -3116

# Row 112 (G=27960)
$ws.Range("H112").Value = 5500
$ws.Range("J112").Value = 5600
$ws.Range("L112").Value = 16800
$ws.Range("N112").Value = -19016

# Row 116 (G=27778)
$ws.Range("H116").Value = 5297.9287
$ws.Range("J116").Value = 5992.3
$ws.Range("L116").Value = 5992.3
$ws.Range("N116").Value = -12876.3

# Row 122 (G=36237)
$ws.Range("H122").Value = 746.43335
$ws.Range("I122").Value = 768.0345
$ws.Range("J122").Value = 120
$ws.Range("K122").Value = 2304.1035
$ws.Range("L122").Value = 360
$ws.Range("M122").Value = 145.8964999999998
$ws.Range("N122").Value = -5260

# Row 125 (G=36228)
$ws.Range("H125").Value = 2149.5
$ws.Range("I125").Value = 2149.5
$ws.Range("K125").Value = 19345.5
$ws.Range("M125").Value = -16885.5

# Row 127 (G=36114)
$ws.Range("H127").Value = 1738.6364
$ws.Range("I127").Value = 1579.2
$ws.Range("K127").Value = 4737.6
$ws.Range("M127").Value = 222.3999999999996

# Row 137 (G=44013)
$ws.Range("H137").Value = 2599.8333
$ws.Range("I137").Value = 2519.8
$ws.Range("K137").Value = 7559.400000000001
$ws.Range("M137").Value = -5009.400000000001

# Row 138 (G=44169)
$ws.Range("H138").Value = 8715.809999999999
$ws.Range("I138").Value = 2023.25
$ws.Range("J138").Value = 10290.529
$ws.Range("K138").Value = 6069.75
$ws.Range("L138").Value = 30871.587
$ws.Range("M138").Value = -929.75
$ws.Range("N138").Value = -41151.587

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (G=27714)
$ws.Range("H45").Value = 2495.8572
$ws.Range("I45").Value = 2263.6924
$ws.Range("J45").Value = 5514
$ws.Range("K45").Value = 2263.6924
$ws.Range("L45").Value = 5514
$ws.Range("M45").Value = -1886.6924
$ws.Range("N45").Value = -6268

# Row 102 (G=19945)
$ws.Range("H102").Value = 1033.4615
$ws.Range("I102").Value = 1037.0834
$ws.Range("K102").Value = 1037.0834
$ws.Range("M102").Value = 584.9166

# Row 122 (G=36168)
$ws.Range("H122").Value = 1130.4667
$ws.Range("I122").Value = 1068.3572
$ws.Range("K122").Value = 3205.0716
$ws.Range("M122").Value = -755.0715999999998

# Row 132 (G=43997)
$ws.Range("H132").Value = 2824
$ws.Range("I132").Value = 1550.8
$ws.Range("J132").Value = 4642.857
$ws.Range("K132").Value = 4652.4
$ws.Range("L132").Value = 13928.571
$ws.Range("M132").Value = -2122.4
$ws.Range("N132").Value = -18988.571

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (G=27706)
$ws.Range("H107").Value = 1750.8572
$ws.Range("I107").Value = 1638.8
$ws.Range("J107").Value = 2031
$ws.Range("K107").Value = 1638.8
$ws.Range("L107").Value = 2031
$ws.Range("M107").Value = 281.2
$ws.Range("N107").Value = -5871

# Row 125 (G=34235)
$ws.Range("H125").Value = 99990
$ws.Range("J125").Value = 99990
$ws.Range("L125").Value = 99990
$ws.Range("N125").Value = -109830

# Row 134 (G=43998)
$ws.Range("H134").Value = 1555
$ws.Range("I134").Value = 1555
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4665
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -2130

$ws = $wb.Worksheets.Item("CRP")
# Row 37 (G=2021)
$ws.Range("H37").Value = 15000
$ws.Range("J37").Value = 15000
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15214

# Row 88 (G=10608)
$ws.Range("H88").Value = 12890
$ws.Range("J88").Value = 12890
$ws.Range("L88").Value = 12890
$ws.Range("N88").Value = -13702

# Row 91 (G=10608)
$ws.Range("H91").Value = 12890
$ws.Range("J91").Value = 12890
$ws.Range("L91").Value = 12890
$ws.Range("N91").Value = -15698

# Row 99 (G=36198)
$ws.Range("H99").Value = 3500.125
$ws.Range("J99").Value = 3066.6667
$ws.Range("L99").Value = 3066.6667
$ws.Range("N99").Value = -6062.6667

# Row 126 (G=36198)
$ws.Range("H126").Value = 3500.125
$ws.Range("J126").Value = 3066.6667
$ws.Range("L126").Value = 9200.000100000001
$ws.Range("N126").Value = -14140.0001

# Row 132 (G=44019)
$ws.Range("H132").Value = 2644.6
$ws.Range("I132").Value = 2597
$ws.Range("K132").Value = 7791
$ws.Range("M132").Value = -5261

# Row 134 (G=44020)
$ws.Range("H134").Value = 6587.3335
$ws.Range("I134").Value = 6587.3335
$ws.Range("K134").Value = 19762.0005
$ws.Range("M134").Value = -17227.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 97 (G=19846)
$ws.Range("H97").Value = 1676.4
$ws.Range("I97").Value = 1699
$ws.Range("J97").Value = 1642.5
$ws.Range("K97").Value = 5097
$ws.Range("L97").Value = 4927.5
$ws.Range("M97").Value = -4601
$ws.Range("N97").Value = -5919.5

# Row 98 (G=19843)
$ws.Range("H98").Value = 646.375
$ws.Range("I98").Value = 357.66666
$ws.Range("J98").Value = 819.6
$ws.Range("K98").Value = 1072.99998
$ws.Range("L98").Value = 2458.8
$ws.Range("M98").Value = 425.0000199999999
$ws.Range("N98").Value = -5454.8

# Row 140 (G=44097)
$ws.Range("H140").Value = 836566.8
$ws.Range("I140").Value = 836566.8
$ws.Range("K140").Value = 2509700.4
$ws.Range("M140").Value = -2504520.4

$ws = $wb.Worksheets.Item("GSM")
# Row 21 (G=4430)
$ws.Range("H21").Value = 8500002
$ws.Range("I21").Value = 10000000
$ws.Range("J21").Value = 7000003.5
$ws.Range("K21").Value = 10000000
$ws.Range("L21").Value = 7000003.5
$ws.Range("M21").Value = -9999827
$ws.Range("N21").Value = -7000349.5

# Row 30 (G=4430)
$ws.Range("H30").Value = 8500002
$ws.Range("I30").Value = 10000000
$ws.Range("J30").Value = 7000003.5
$ws.Range("K30").Value = 10000000
$ws.Range("L30").Value = 7000003.5
$ws.Range("M30").Value = -9999895
$ws.Range("N30").Value = -7000213.5

# Row 126 (G=36184)
$ws.Range("H126").Value = 4482.75
$ws.Range("I126").Value = 4032.4
$ws.Range("J126").Value = 5233.3335
$ws.Range("K126").Value = 12097.2
$ws.Range("L126").Value = 15700.0005
$ws.Range("M126").Value = -9627.200000000001
$ws.Range("N126").Value = -20640.0005

# Row 132 (G=44008)
$ws.Range("H132").Value = 2643.3333
$ws.Range("I132").Value = 2643.3333
$ws.Range("K132").Value = 7929.999899999999
$ws.Range("M132").Value = -5399.999899999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (G=36249)
$ws.Range("H7").Value = 1749.25
$ws.Range("I7").Value = 1749.25
$ws.Range("K7").Value = 1749.25
$ws.Range("M7").Value = -1637.25

# Row 21 (G=2672)
$ws.Range("H21").Value = 4255.25
$ws.Range("I21").Value = 4499.5
$ws.Range("J21").Value = 4011
$ws.Range("K21").Value = 4499.5
$ws.Range("L21").Value = 4011
$ws.Range("M21").Value = -4325.5
$ws.Range("N21").Value = -4359

# Row 33 (G=4106)
$ws.Range("H33").Value = 1300
$ws.Range("I33").Value = 1300
$ws.Range("K33").Value = 1300
$ws.Range("M33").Value = -1010

# Row 61 (G=27740)
$ws.Range("H61").Value = 2130.5
$ws.Range("I61").Value = 2145.5557
$ws.Range("J61").Value = 1995
$ws.Range("K61").Value = 2145.5557
$ws.Range("L61").Value = 1995
$ws.Range("M61").Value = -1943.5557
$ws.Range("N61").Value = -2399

# Row 113 (G=27740)
$ws.Range("H113").Value = 2130.5
$ws.Range("I113").Value = 2145.5557
$ws.Range("J113").Value = 1995
$ws.Range("K113").Value = 2145.5557
$ws.Range("L113").Value = 1995
$ws.Range("M113").Value = 24.44430000000011
$ws.Range("N113").Value = -6335

# Row 125 (G=34271)
$ws.Range("H125").Value = 56000
$ws.Range("J125").Value = 56000
$ws.Range("L125").Value = 56000
$ws.Range("N125").Value = -65840

# Row 126 (G=36249)
$ws.Range("H126").Value = 1749.25
$ws.Range("I126").Value = 1749.25
$ws.Range("K126").Value = 5247.75
$ws.Range("M126").Value = -2777.75
